$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) text in place.
$ws.Range("B1").Value = "序號"
$ws.Range("C1").Value = "候選人"
$ws.Range("D1").Value = "得票數"
$ws.Range("E1").ClearContents()
$ws.Range("F1").ClearContents()
$ws.Range("F1").Style = "一般"

# Move the active selection to C2 (reflecting the last user click before save).
$ws.Range("C2").Select()
